# B6-PowerPoint.pptx edit, 21 May 2020
#
# 1) Three tables (slides 14-16) get re-styled from the old "Table_0"
#    style ({3E4481D4-...}) to the "Best match for document" style
#    ({1B6C281B-...}) that PowerPoint wrote when the table style was
#    changed from the Table Design gallery.
# 2) The deck's main theme (slideMaster1 -> theme1.xml, the "Integral"/
#    "Red Violet" theme) was swapped for the stock Office theme - i.e.
#    the 12-slot theme colour scheme now matches the default "Office"
#    palette (the font scheme/format scheme were already identical
#    between the two themes, so only the colours actually change).

$p = $ppt.ActivePresentation

# --- 1. Re-apply table style on the three affected tables -------------
$targetStyle = "{1B6C281B-59E3-4059-97E9-88E056D7EBD5}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($targetStyle)
        }
    }
}

# --- 2. Swap the main theme's colour scheme for the Office defaults ---
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

# index -> (slot, new RGB as VBA-style BGR-int = R + G*256 + B*65536)
$themeColors.Colors(1).RGB  = 0         # dk1      000000
$themeColors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388   # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407     # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456   # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797  # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477   # folHlink 954F72
